# Add a new column H ("d_oilp_u") to Sheet1.
# H1 gets the same header style as the existing header cells (A1:G1);
# H121:H278 get the computed "unexpected oil-price change" series.
# (Rows 2-120 have no value for this series in the source data.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from G1 onto H1, then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "d_oilp_u"

# Populate H121:H278 with the new series values.
$values = @(-0.05382291278422713, -0.01669820730076221, -0.05154042153572558, -0.07330281913115932, -0.03630920940450899, -0.1501305285132455, -0.02379360905419148, -0.08356545131594384, 0.07415175403240948, 0.03783463392486475, -0.0363449083722216, -0.01781803690666539, 0.1087355284495679, -0.0424035323192955, -0.02625633456655407, 0.05227542720551348, -0.04362984021265959, -0.0917186834974828, 0.003057979985370274, 0.02854698088935237, 0.1503535656022139, 0.07821306196594424, -0.1257674410424183, -0.07051458767246466, 0.03418065587437091, 0.003739003945602448, -0.007645436349183576, 0.006275885804139847, -0.05528622443033893, 0.05793571663850994, 0.01613905090508005, 0.0590743644361762, 0.0140043009790638, -0.02293950697215674, -0.0592090956776774, -0.06448669558237885, 0.07877662975573152, 0.00487166037504938, -0.03904731778250259, 0.01795282190056557, 0.01332982802975735, -0.007612365018409228, 0.0279463048087365, -0.02377596374770885, 0.00811231079481356, -0.02286184072087583, 0.03957216582793954, 0.04547787216883137, 0.07927915879421299, 0.1068955826387246, 0.1545308128784022, 0.2899725958325545, 0.1149323666230981, -0.239225221060106, 0.1428084032233237, -0.1850511106516688, 0.01179672535759213, 0.03425551048742026, 0.1505064816581361, 0.09652619663096829, 0.02429965815209867, -0.009963261971257076, 0.07889394246966264, 0.1984964439455705, 0.1039078991141467, -0.06947660846966297, -0.05831889548829139, -0.2080905265273816, -0.07325381533985009, 0.0181900108435844, 0.1703972042611319, -0.114656580788099, -0.04479513379283206, 0.03268193049942525, -0.04362742017962606, -0.1319827020263418, -0.009123980993370751, 0.03087013432162378, 0.03323598108910009, 0.0248111966962874, 0.02481924135299796, 0.04394212185649859, -0.1158521603702471, 0.0009668375499352067, -0.0848971942372394, -0.08316832868799384, -0.03502928839515818, -0.05229547638717591, 0.0008995502855464466, 0.04353177736200298, -0.05841965097616875, -0.07814255065180387, -0.03465294349354053, -0.00906767498259331, 0.04933976572215037, -0.05520947539416721, -0.06519671259753146, 0.1131483677038618, 0.2355340927000604, 0.1345433261207205, -0.2181528489921698, -0.02772516844827155, -0.05850498436312446, -0.06296180157175701, 0.07690084319433232, 0.02105340919783316, 0.02104377746705755, 0.04410858325627576, 0.03524088348828425, -0.003497380532170524, -0.07985086115016582, -0.01704500252928121, 0.1269352468844351, 0.1778429933861054, 1.173602289544946, -0.217435697495469, -0.6356128268659447, -0.2194141013227782, 0.01345056183948889, -0.0345134063792818, 0.1078272164645631, 0.08689912212606998, -0.2342881158672752, -0.05715841383994835, -0.02563478672672703, -0.1021180662437211, -0.1340848266969799, 0.06727916974911086, -0.1257965014804521, -0.005391634190396566, -0.08159539027204854, 0.04422652431955587, 0.03925991313760679, -0.1100524615415281, -0.08053420429059965, 0.1637230235625795, -0.08240593822990672, -0.1662906875089201, -0.08290398854088821, 0.1832783680392476, 0.04690998243445321, -0.1911833023324556, 0.05434347979122833, 0.1455869840503832, 0.04465145363767142, 0.03310059793916142, -0.06905803717072523, 0.09968569698922103, 0.1634086658132228, 0.0486375817110325, -0.0639487246002739, -0.07647855894149558, 0.03863310547273358, 0.0920518284066576, -0.03260215914986198, -0.1808702396825916, -0.04223918717973874, -0.05212559875880451)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(121 + $i, 8).Value = $values[$i]
}
